$d = $word.ActiveDocument

# --- Paragraph 3 (1-based): sampling frequency bullet ---
$p3xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="004A595D" w:rsidRDefault="00D116DA" w:rsidP="001644A3"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>I chose a sampling frequency of 1</w:t></w:r><w:r><w:t>6</w:t></w:r><w:r><w:t>kHz</w:t></w:r><w:r><w:t xml:space="preserve"> because it</w:t></w:r><w:r><w:t xml:space="preserve"> allows the PSD to</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>contain several sidelobes of the 4kHz modulated signal.</w:t></w:r></w:p>'
$d.Paragraphs.Item(3).Range.InsertXML($p3xml)

# --- Paragraph 4 (1-based): band-limit bullet ---
$p4xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="001644A3" w:rsidRDefault="001644A3" w:rsidP="001644A3"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">To ensure the PSD was averaged over multiple data blocks, I </w:t></w:r><w:r><w:t xml:space="preserve">have a variable to set the number of symbols per PSD block. Adjusting this variable has the effect of visually varying the noise </w:t></w:r><w:r><w:t xml:space="preserve">present </w:t></w:r><w:r><w:t>in the PSD due to the random bits.</w:t></w:r></w:p>'
$d.Paragraphs.Item(4).Range.InsertXML($p4xml)

# --- Paragraph 5 (1-based): en-dash -> hyphen bullet ---
$p5xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="004A62C6" w:rsidRDefault="002E7640" w:rsidP="001644A3"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>-</w:t></w:r></w:p>'
$d.Paragraphs.Item(5).Range.InsertXML($p5xml)

# --- Paragraph 6 (1-based): main difference bullet, full rewrite ---
$p6xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="002E7640" w:rsidRDefault="002E7640" w:rsidP="001644A3"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>The main difference I’m noticing between s, z1 and z2 is that z1 and z2 are ‘rectangular pulses’ with a lower number of harmonics than the ideally generated pulse train p(t-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>).</w:t></w:r><w:r><w:t xml:space="preserve"> For z2, this is because its real and imaginary parts u(t) and v(t) are LP filtered as a part of the quadrature demodulation, which both cuts out the negative signal and some of the harmonics of each signal. My LPF had a cutoff of 2kHz, which only retained one sidelobe, visualized in figure 2 b. For z1, this is because finding the analytical signal x1+ of the received signal x1 neutralizes any negative frequency components, or any components that were aliased effectively </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bandlimiting</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> the signal when it is shifted back down to baseband, as seen in Figure 1 c/d. This </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bandlimiting</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> can be asymmetrical depending on the sampling frequency and the carrier frequency. In this instance, I chose a sampling frequency of 16kHz which causes the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bandlimiting</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> effect to be symmetrical. </w:t></w:r><w:r><w:t xml:space="preserve"> In this case, the PSD of z1 was ‘bandlimited’ to 4 sidelobes, which makes the </w:t></w:r><w:r><w:t>signal z1(t) appear to be a square pulse with more harmonics than signal z2(t).</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$d.Paragraphs.Item(6).Range.InsertXML($p6xml)
